# Applies the "SDKv5.00" summary update to the tinyyolov3_cam summary sheet.
# - S6 and S32 (timing figures) are updated to new measured values.
# - B8..B26 (layer "name" column) swap the long ONNX node-path strings for
#   short numeric-looking labels (kept as text, not numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated timing values (column S) ---
$ws.Range("S6").Value = 1664
$ws.Range("S32").Value = 28307

# --- Updated "name" column values (column B), stored as text ---
$nameUpdates = @{
    8  = "146"
    9  = "78"
    10 = "149"
    11 = "82"
    12 = "152"
    13 = "86"
    14 = "155"
    15 = "90"
    16 = "158"
    17 = "94"
    18 = "161"
    19 = "121"
    20 = "164"
    21 = "167"
    22 = "170"
    24 = "173"
    25 = "140"
    26 = "176"
}

foreach ($row in $nameUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 2)
    # Force text storage (values look numeric but must remain strings),
    # then restore the default "Normal" style so no new cell formatting
    # is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $nameUpdates[$row]
    $cell.Style = "Normal"
}

# --- Page margins reset to Excel's standard defaults (inches -> points) ---
$ws.PageSetup.LeftMargin = 0.75 * 72
$ws.PageSetup.RightMargin = 0.75 * 72
$ws.PageSetup.TopMargin = 1 * 72
$ws.PageSetup.BottomMargin = 1 * 72
$ws.PageSetup.HeaderMargin = 0.5 * 72
$ws.PageSetup.FooterMargin = 0.5 * 72
